$p = $ppt.ActivePresentation

# Slide 6: "Project Infrastructure & Software" / "Presented by Brian Bauer"
$s6 = $p.Slides.Add(6, 2)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Project Infrastructure & Software"
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "Presented by Brian Bauer"

# Slide 7: "Graphics and Visual Modeling" / "Presented by David Kalbfleisch"
$s7 = $p.Slides.Add(7, 2)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Graphics and Visual Modeling"
$tr7 = $s7.Shapes.Item(2).TextFrame.TextRange
$tr7.Text = "Presented by David "
[void]$tr7.InsertAfter("Kalbfleisch")

# Slide 8: "Neural Network Development and Training" / "Presented by Iman Ismail"
$s8 = $p.Slides.Add(8, 2)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Neural Network Development and Training"
$s8.Shapes.Item(2).TextFrame.TextRange.Text = "Presented by Iman Ismail"

# Slide 9: "Conclusion and Final Remarks" / "Demo to follow"
$s9 = $p.Slides.Add(9, 2)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = "Conclusion and Final Remarks"
$tr9 = $s9.Shapes.Item(2).TextFrame.TextRange
$tr9.Text = "Demo "
[void]$tr9.InsertAfter("to follow")

[void]$p.Save()
